# Merged jacobs files in - fills in grades/examples for the "Peer  and self
# assessment" sheet (the active/tab-selected sheet) and updates the active
# selection, matching the authoritative OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Row 2 - Self assessment / Criterion 1 Online collaboration
$ws.Range("B2").Value = "Good"
$ws.Range("C2").Value = "I have been active in the weekly online Discord meetings - both in `nterms of planning (Adding discussion points to the meeting agenda, taking responsibility for leading a group discussion or a walkthrough of conducted work, etc) `nI have been reading through literature and materials uploaded`nto GiT by other project group members, and actively been adding comments and suggestions. I have been uploading relevant materials to`nGiT as well.  `nI have been responding fairly quickly to messages, both private and group announcements, on the discord platform. "

# Row 3 - Peer assessment 9
$ws.Range("B3").Value = "Good"
$ws.Range("C3").Value = "Research, hardware setup "

# Row 15 - Self assessment (final reflection section)
$ws.Range("B15").Value = "Good"
$ws.Range("C15").Value = "Active collaboration with Ahmet about Machine learning (Including sharing  research and literature), Active collaboration with Alex and Morcel about the PoA (Including giving inputs and suggestions for conducted work), Active communication and use of GiT and Discord, "

# Restore the saved window scroll position (topLeftCell = B10) and selection
# (activeCell = C13) recorded for this sheet view.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("C13").Select() | Out-Null
